{"js": "// Fix report typos and merge accidentally split runs.\n//\n// 1. \"Which are the most frequent words found in the tweets based on the\n//    top 10 [dis][approvals?]\" -> the \"dis\" + \"approvals?\" runs are merged\n//    back into a single run (no character changes, just a run merge).\n// 2. \"What is the correlation between the number of retweets and the\n//    [dis][approval ratings for the top 10 approvals?]\" -> same kind of\n//    run merge.\n// 3. \"MOST FREQUANT\" -> \"MOST FREQUENT\" (typo fix, 4 occurrences).\n// 4. \"APPORVALS\" -> \"APPROVALS\" and \"DISAPPORVALS\" -> \"DISAPPROVALS\"\n//    (typo fix, 2 occurrences each; highlight formatting is preserved\n//    because the replace happens in-place on the matched range).\n\nconst body = context.document.body;\n\nasync function replaceAll(searchText, replacement, options) {\n  const searchOptions = Object.assign({ matchCase: true, matchWholeWord: false }, options || {});\n  const results = body.search(searchText, searchOptions);\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(replacement, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\n// --- Typo fixes -----------------------------------------------------------\nawait replaceAll(\"MOST FREQUANT WORDS\", \"MOST FREQUENT WORDS\");\n// \"APPORVALS\" also matches the \"APPORVALS\" inside \"DISAPPORVALS\", so this\n// single pass fixes both APPORVALS -> APPROVALS and DISAPPORVALS ->\n// DISAPPROVALS while keeping each run's own highlight color.\nawait replaceAll(\"APPORVALS\", \"APPROVALS\");\n\n// --- Merge the two paragraphs whose \"dis\" prefix was split into its own\n//     run back into a single contiguous run -------------------------------\nawait replaceAll(\n  \"Which are the most frequent words found in the tweets based on the top 10 disapprovals?\",\n  \"Which are the most frequent words found in the tweets based on the top 10 disapprovals?\"\n);\n\nawait replaceAll(\n  \"What is the correlation between the number of retweets and the disapproval ratings for the top 10 approvals?\",\n  \"What is the correlation between the number of retweets and the disapproval ratings for the top 10 approvals?\"\n);\n", "ps1": "# Fix report typos and merge accidentally split runs.\n#\n# 1. \"Which are the most frequent words found in the tweets based on the\n#    top 10 [dis][approvals?]\" -> the \"dis\" + \"approvals?\" runs are merged\n#    back into a single run (no character changes, just a run merge).\n# 2. \"What is the correlation between the number of retweets and the\n#    [dis][approval ratings for the top 10 approvals?]\" -> same kind of\n#    run merge.\n# 3. \"MOST FREQUANT\" -> \"MOST FREQUENT\" (typo fix, 4 occurrences).\n# 4. \"APPORVALS\" -> \"APPROVALS\" and \"DISAPPORVALS\" -> \"DISAPPROVALS\"\n#    (typo fix, 2 occurrences each; highlight formatting is preserved\n#    because the replace happens in-place on the matched range).\n\n$d = $word.ActiveDocument\n\nfunction Replace-All($findText, $replaceText) {\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.Text = $replaceText\n    $find.Forward = $true\n    $find.Wrap = 1              # wdFindContinue\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2) | Out-Null\n}\n\n# --- Typo fixes -------------------------------------------------------\nReplace-All \"MOST FREQUANT WORDS\" \"MOST FREQUENT WORDS\"\n# \"APPORVALS\" also matches the \"APPORVALS\" inside \"DISAPPORVALS\", so this\n# single pass fixes both APPORVALS -> APPROVALS and DISAPPORVALS ->\n# DISAPPROVALS while keeping each run's own highlight color.\nReplace-All \"APPORVALS\" \"APPROVALS\"\n\n# --- Merge the two paragraphs whose \"dis\" prefix was split into its own\n#     run back into a single contiguous run ----------------------------\nReplace-All \"Which are the most frequent words found in the tweets based on the top 10 disapprovals?\" \"Which are the most frequent words found in the tweets based on the top 10 disapprovals?\"\nReplace-All \"What is the correlation between the number of retweets and the disapproval ratings for the top 10 approvals?\" \"What is the correlation between the number of retweets and the disapproval ratings for the top 10 approvals?\"\n"}
